{"js": "// The site rebuild dropped the trailing \"\u00a9 2020 ... Creative Commons\n// Attribution\" footer block (and the blank/page-break paragraphs that\n// introduced it) that used to be duplicated right after the\n// bibliography's last entry (\"...McGraw-Hill, 1990.\"). An identical\n// blank + page-break + nothing-following footer sequence still exists\n// at the very end of the document and must be left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load the text of every paragraph so we can locate the copyright\n// notice paragraph unambiguously.\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst marker = \"\\u00A9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === marker) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the copyright footer paragraph to remove.\");\n}\n\n// The two paragraphs immediately preceding it are the blank paragraph\n// and the blank page-break paragraph that belong to this same (now\n// removed) footer block.\nconst blankPageBreakPara = target.getPrevious();\nconst blankPara = blankPageBreakPara.getPrevious();\n\ntarget.delete();\nblankPageBreakPara.delete();\nblankPara.delete();\n\nawait context.sync();\n", "ps1": "# The site rebuild dropped the trailing \"(c) 2020 ... Creative Commons\n# Attribution\" footer block (and the blank / page-break paragraphs that\n# introduced it) that used to be duplicated right after the\n# bibliography's last entry (\"...McGraw-Hill, 1990.\"). An identical\n# blank + page-break + footer-less sequence still exists at the very\n# end of the document and must be left untouched.\n\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Contact: luizeleno@usp.br\")\nif (-not $found) {\n    throw \"Could not locate the copyright footer paragraph to remove.\"\n}\n$targetStart = $findRange.Start\n$targetEnd = $findRange.End\n\n# Map the Find hit back to its paragraph index in $d.Paragraphs (more\n# reliable here than a sub-range's own .Paragraphs collection).\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -le $targetStart -and $p.Range.End -ge $targetEnd) {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not map the found text back to a paragraph.\"\n}\n\n# Delete the copyright paragraph and the two blank paragraphs right\n# before it (the blank page-break paragraph, then the plain blank\n# paragraph). Go from the highest index down so earlier indices stay\n# valid as each delete happens.\n$d.Paragraphs.Item($targetIndex).Range.Delete()\n$d.Paragraphs.Item($targetIndex - 1).Range.Delete()\n$d.Paragraphs.Item($targetIndex - 2).Range.Delete()\n"}
